$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.486.51'
$ws.Range("E2").Value = '  -1.11%  '

$ws.Range("D3").Value = '3.082.11'
$ws.Range("E3").Value = '  +0.38%  '

$ws.Range("E4").Value = '  +0.09%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '555.19'
$ws.Range("E5").Value = '  +0.88%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '137.22'
$ws.Range("E6").Value = '  -2.92%  '

$ws.Range("E7").Value = '  +0.02%  '

$ws.Range("D8").Value = '3.075.77'
$ws.Range("E8").Value = '  +0.65%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.493'
$ws.Range("E9").Value = '  +0.79%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.70'
$ws.Range("E10").Value = '  +4.03%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.154'
$ws.Range("E11").Value = '  +0.51%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.451'
$ws.Range("E12").Value = '  -0.86%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '35.17'
$ws.Range("E13").Value = '  -1.55%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000215'
$ws.Range("E14").Value = '  -1.10%  '

$ws.Range("D15").Value = '3.586.94'
$ws.Range("E15").Value = '  +0.55%  '

$ws.Range("D16").Value = '63.598.46'
$ws.Range("E16").Value = '  -1.06%  '

$ws.Range("E17").Value = '  +0.00%  '

$ws.Range("D18").Value = '3.101.41'
$ws.Range("E18").Value = '  +0.76%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '503.65'
$ws.Range("E19").Value = '  +3.09%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.61'
$ws.Range("E20").Value = '  -0.18%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.58'
$ws.Range("E21").Value = '  -0.44%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.702'
$ws.Range("E22").Value = '  +2.29%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.20'
$ws.Range("E23").Value = '  -0.63%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.31'
$ws.Range("E24").Value = '  -1.00%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '77.22'
$ws.Range("E25").Value = '  -1.32%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.998'
$ws.Range("E26").Value = '  -0.12%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.76'
$ws.Range("E27").Value = '  +2.00%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.27'
$ws.Range("E28").Value = '  +1.61%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.03'
$ws.Range("E29").Value = '  -2.52%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.00'
$ws.Range("E30").Value = '  +0.08%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '26.00'
$ws.Range("E31").Value = '  +0.42%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.53'
$ws.Range("E32").Value = '  -4.04%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.10'
$ws.Range("E33").Value = '  -1.99%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '532.02'
$ws.Range("E34").Value = '  -10.52%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '56.63'
$ws.Range("E35").Value = '  +8.52%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.87'
$ws.Range("E36").Value = '  -2.73%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.15'
$ws.Range("E37").Value = '  -5.20%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0412'
$ws.Range("E38").Value = '  +2.37%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0794'
$ws.Range("E39").Value = '  +0.25%  '

$ws.Range("D40").Value = '3.057.84'
$ws.Range("E40").Value = '  +2.58%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.117'
$ws.Range("E41").Value = '  -1.09%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.07'
$ws.Range("E42").Value = '  -1.80%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.57'
$ws.Range("E43").Value = '  -10.52%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.252'
$ws.Range("E44").Value = '  +3.22%  '

$ws.Range("E45").Value = '  +0.05%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.07'
$ws.Range("E46").Value = '  -2.09%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '121.30'
$ws.Range("E47").Value = '  +1.46%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '23.94'
$ws.Range("E48").Value = '  -5.53%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.106'
$ws.Range("E49").Value = '  -1.23%  '

$ws.Range("D50").Value = '0.0₃0495'
$ws.Range("E50").Value = '  -7.01%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.01'
$ws.Range("E51").Value = '  -3.32%  '
